$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Business partner check / Employee module related field updates on row 2
$ws.Range("C2").Value = "Camera"
$ws.Range("S2").Value = "Camera"
$ws.Range("T2").Value = "Troy"
$ws.Range("U2").Value = "KKThoppu"
$ws.Range("AD2").Value = "Sales"

# Update the active selection recorded in the sheet view
$ws.Range("AG2").Select() | Out-Null
